# Auto-generated: apply cryptos list update (commit: "Updated cryptos list on Sun Apr 23 06:42:47 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the source file's inline-string cells),
# even when the text looks like a number (e.g. "1.004", "0.00001047", thousand-dot
# separated prices like "27.688.94"). Without this, Excel auto-types such strings
# as numbers. NumberFormat is forced to Text first, then reset to the default cell
# style afterwards so no stray formatting is left behind on the cell.
function Set-TextValue($sheet, $addr, $val) {
    $cell = $sheet.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "27.688.94"
$ws.Range("E2").Value = "  +1.21%  "
# Row 3
$ws.Range("D3").Value = "1.874.80"
$ws.Range("E3").Value = "  +0.89%  "
# Row 4
Set-TextValue $ws "D4" "1.004"
$ws.Range("E4").Value = "  +0.21%  "
# Row 5
Set-TextValue $ws "D5" "332.15"
$ws.Range("E5").Value = "  +2.69%  "
# Row 6
$ws.Range("E6").Value = "  +0.17%  "
# Row 7
Set-TextValue $ws "D7" "0.4723"
$ws.Range("E7").Value = "  +4.23%  "
# Row 8
Set-TextValue $ws "D8" "0.3950"
$ws.Range("E8").Value = "  +2.23%  "
# Row 9
Set-TextValue $ws "D9" "47.87"
$ws.Range("E9").Value = "  -1.08%  "
# Row 10
Set-TextValue $ws "D10" "0.08044"
$ws.Range("E10").Value = "  +1.66%  "
# Row 11
Set-TextValue $ws "D11" "1.027"
$ws.Range("E11").Value = "  +1.02%  "
# Row 12
Set-TextValue $ws "D12" "22.05"
$ws.Range("E12").Value = "  +3.18%  "
# Row 13
$ws.Range("D13").Value = "1.883.79"
$ws.Range("E13").Value = "  +0.81%  "
# Row 14
Set-TextValue $ws "D14" "5.966"
$ws.Range("E14").Value = "  +1.09%  "
# Row 15
Set-TextValue $ws "D15" "7.141"
$ws.Range("E15").Value = "  +0.31%  "
# Row 16
$ws.Range("E16").Value = "  +0.04%  "
# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws "D17" "0.00001047"
$ws.Range("E17").Value = "  +1.50%  "
# Row 18
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D18" "87.06"
$ws.Range("E18").Value = "  +1.50%  "
# Row 19
Set-TextValue $ws "D19" "0.06665"
$ws.Range("E19").Value = "  +2.29%  "
# Row 20
$ws.Range("E20").Value = "  +0.80%  "
# Row 21
$ws.Range("E21").Value = "  +0.17%  "
# Row 22
$ws.Range("D22").Value = "27.707.09"
$ws.Range("E22").Value = "  +1.26%  "
# Row 23
Set-TextValue $ws "D23" "5.521"
$ws.Range("E23").Value = "  -0.12%  "
# Row 24
Set-TextValue $ws "D24" "10.99"
$ws.Range("E24").Value = "  +1.39%  "
# Row 25
Set-TextValue $ws "D25" "2.306"
$ws.Range("E25").Value = "  +1.13%  "
# Row 26
$ws.Range("D26").Value = "2.105.28"
$ws.Range("E26").Value = "  +0.89%  "
# Row 27
Set-TextValue $ws "D27" "158.53"
$ws.Range("E27").Value = "  +3.02%  "
# Row 28
Set-TextValue $ws "D28" "20.21"
$ws.Range("E28").Value = "  +2.40%  "
# Row 29
Set-TextValue $ws "D29" "2.100"
$ws.Range("E29").Value = "  +1.36%  "
# Row 30
Set-TextValue $ws "D30" "5.591"
$ws.Range("E30").Value = "  +2.94%  "
# Row 31
Set-TextValue $ws "D31" "122.21"
$ws.Range("E31").Value = "  +1.29%  "
# Row 32
Set-TextValue $ws "D32" "0.9742"
$ws.Range("E32").Value = "  +4.29%  "
# Row 33
Set-TextValue $ws "D33" "0.09558"
$ws.Range("E33").Value = "  +2.93%  "
# Row 34
Set-TextValue $ws "D34" "1.449"
$ws.Range("E34").Value = "  -2.56%  "
# Row 35
Set-TextValue $ws "D35" "3.594"
$ws.Range("E35").Value = "  +0.02%  "
# Row 36
Set-TextValue $ws "D36" "5.337"
$ws.Range("E36").Value = "  +1.58%  "
# Row 37
Set-TextValue $ws "D37" "0.06103"
$ws.Range("E37").Value = "  +1.99%  "
# Row 39
Set-TextValue $ws "D39" "1.233"
$ws.Range("E39").Value = "  +0.84%  "
# Row 40
Set-TextValue $ws "D40" "8.257"
$ws.Range("E40").Value = "  +0.52%  "
# Row 41
Set-TextValue $ws "D41" "0.6023"
$ws.Range("E41").Value = "  +2.04%  "
# Row 42
Set-TextValue $ws "D42" "0.1906"
$ws.Range("E42").Value = "  +1.07%  "
# Row 43
Set-TextValue $ws "D43" "10.25"
$ws.Range("E43").Value = "  +1.59%  "
# Row 44
$ws.Range("E44").Value = "  -1.62%  "
# Row 45
Set-TextValue $ws "D45" "0.5717"
$ws.Range("E45").Value = "  +1.77%  "
# Row 46
Set-TextValue $ws "D46" "12.15"
$ws.Range("E46").Value = "  +1.89%  "
# Row 47
Set-TextValue $ws "D47" "1.948"
$ws.Range("E47").Value = "  +1.27%  "
# Row 48
Set-TextValue $ws "D48" "3.389"
$ws.Range("E48").Value = "  +0.43%  "
# Row 49
Set-TextValue $ws "D49" "115.71"
$ws.Range("E49").Value = "  +7.03%  "
# Row 50
Set-TextValue $ws "D50" "0.06881"
$ws.Range("E50").Value = "  +1.58%  "
# Row 51
Set-TextValue $ws "D51" "0.00000000303"
$ws.Range("E51").Value = "  +15.73%  "
